$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 749.6667
$ws.Range("I12").Value = 524.6667
$ws.Range("J12").Value = 1199.6666
$ws.Range("K12").Value = 524.6667
$ws.Range("L12").Value = 1199.6666
$ws.Range("M12").Value = -354.6667
$ws.Range("N12").Value = -1539.6666
$ws.Range("H53").Value = 12092.3
$ws.Range("J53").Value = 105
$ws.Range("L53").Value = 105
$ws.Range("N53").Value = -1379
$ws.Range("H64").Value = 3230
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3230
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("H98").Value = 6621.3335
$ws.Range("I98").Value = 6682
$ws.Range("K98").Value = 6682
$ws.Range("M98").Value = -5184
$ws.Range("H111").Value = 669.2
$ws.Range("I111").Value = 669.2
$ws.Range("K111").Value = 2007.6
$ws.Range("M111").Value = 1059.4
$ws.Range("H115").Value = 685
$ws.Range("I115").Value = 685
$ws.Range("K115").Value = 2055
$ws.Range("M115").Value = -488
$ws.Range("H122").Value = 6621.3335
$ws.Range("I122").Value = 6682
$ws.Range("K122").Value = 20046
$ws.Range("M122").Value = -17596
$ws.Range("H137").Value = 54081.684
$ws.Range("I137").Value = 1178
$ws.Range("K137").Value = 3534
$ws.Range("M137").Value = -984
$ws.Range("H138").Value = 1614.57
$ws.Range("I138").Value = 1003.74286
$ws.Range("J138").Value = 1943.4769
$ws.Range("K138").Value = 3011.22858
$ws.Range("L138").Value = 5830.4307
$ws.Range("M138").Value = 2128.77142
$ws.Range("N138").Value = -16110.4307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4651.674
$ws.Range("I32").Value = 2888.5059
$ws.Range("J32").Value = 19638.6
$ws.Range("K32").Value = 2888.5059
$ws.Range("L32").Value = 19638.6
$ws.Range("M32").Value = -2601.5059
$ws.Range("N32").Value = -20212.6
$ws.Range("H45").Value = 1352.6522
$ws.Range("I45").Value = 1178.9231
$ws.Range("J45").Value = 1578.5
$ws.Range("K45").Value = 1178.9231
$ws.Range("L45").Value = 1578.5
$ws.Range("M45").Value = -801.9231
$ws.Range("N45").Value = -2332.5
$ws.Range("H61").Value = 25355.854
$ws.Range("I61").Value = 27656.633
$ws.Range("K61").Value = 27656.633
$ws.Range("M61").Value = -27444.633
$ws.Range("H74").Value = 650.1429000000001
$ws.Range("I74").Value = 524.53656
$ws.Range("K74").Value = 524.53656
$ws.Range("M74").Value = 349.46344
$ws.Range("H77").Value = 650.1429000000001
$ws.Range("I77").Value = 524.53656
$ws.Range("K77").Value = 2622.6828
$ws.Range("M77").Value = 1745.3172
$ws.Range("H97").Value = 1072.6923
$ws.Range("I97").Value = 916.9
$ws.Range("K97").Value = 916.9
$ws.Range("M97").Value = -420.9
$ws.Range("H110").Value = 590.3077
$ws.Range("I110").Value = 614.5
$ws.Range("K110").Value = 614.5
$ws.Range("M110").Value = 1430.5
$ws.Range("H132").Value = 1618.9014
$ws.Range("I132").Value = 1307.6666
$ws.Range("K132").Value = 3922.9998
$ws.Range("M132").Value = -1392.9998
$ws.Range("H136").Value = 25355.854
$ws.Range("I136").Value = 27656.633
$ws.Range("K136").Value = 82969.899
$ws.Range("M136").Value = -80419.899

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2736.0312
$ws.Range("I20").Value = 2492.8
$ws.Range("J20").Value = 3141.4167
$ws.Range("K20").Value = 2492.8
$ws.Range("L20").Value = 3141.4167
$ws.Range("M20").Value = -2245.8
$ws.Range("N20").Value = -3635.4167
$ws.Range("H94").Value = 1044.5714
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -5902
$ws.Range("H107").Value = 1064.375
$ws.Range("I107").Value = 1133.5714
$ws.Range("K107").Value = 1133.5714
$ws.Range("M107").Value = 786.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H23").Value = 91675
$ws.Range("J23").Value = 91675
$ws.Range("L23").Value = 91675
$ws.Range("N23").Value = -92155
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H27").Value = 91675
$ws.Range("J27").Value = 91675
$ws.Range("L27").Value = 91675
$ws.Range("N27").Value = -92059
$ws.Range("H31").Value = 1868.1
$ws.Range("I31").Value = 1200
$ws.Range("K31").Value = 1200
$ws.Range("M31").Value = -905
$ws.Range("H34").Value = 1868.1
$ws.Range("I34").Value = 1200
$ws.Range("K34").Value = 1200
$ws.Range("M34").Value = -998
$ws.Range("H62").Value = 2686.4285
$ws.Range("I62").Value = 2179.8
$ws.Range("K62").Value = 2179.8
$ws.Range("M62").Value = -1555.8
$ws.Range("H65").Value = 2686.4285
$ws.Range("I65").Value = 2179.8
$ws.Range("K65").Value = 10899
$ws.Range("M65").Value = -7779
$ws.Range("H107").Value = 889.51514
$ws.Range("I107").Value = 690.64
$ws.Range("K107").Value = 690.64
$ws.Range("M107").Value = 1229.36
$ws.Range("H122").Value = 3937.5
$ws.Range("I122").Value = 2250.1667
$ws.Range("K122").Value = 6750.500100000001
$ws.Range("M122").Value = -4300.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 821785.5600000001
$ws.Range("I4").Value = 1050160
$ws.Range("J4").Value = 250849.5
$ws.Range("K4").Value = 3150480
$ws.Range("L4").Value = 752548.5
$ws.Range("M4").Value = -3150368
$ws.Range("N4").Value = -752772.5
$ws.Range("H122").Value = 1029.9375
$ws.Range("J122").Value = 1129.2307
$ws.Range("L122").Value = 10163.0763
$ws.Range("N122").Value = -15063.0763
$ws.Range("H131").Value = 44181.895
$ws.Range("J131").Value = 55863.867
$ws.Range("L131").Value = 167591.601
$ws.Range("N131").Value = -177671.601

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2124.375
$ws.Range("I102").Value = 2499.1667
$ws.Range("K102").Value = 2499.1667
$ws.Range("M102").Value = -877.1667000000002
$ws.Range("H110").Value = 63867.777
$ws.Range("J110").Value = 63867.777
$ws.Range("L110").Value = 63867.777
$ws.Range("N110").Value = -72047.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2573.6667
$ws.Range("J46").Value = 2794.7144
$ws.Range("L46").Value = 2794.7144
$ws.Range("N46").Value = -3170.7144
$ws.Range("H55").Value = 573.55
$ws.Range("I55").Value = 451.41666
$ws.Range("J55").Value = 756.75
$ws.Range("K55").Value = 451.41666
$ws.Range("L55").Value = 756.75
$ws.Range("M55").Value = -278.41666
$ws.Range("N55").Value = -1102.75
$ws.Range("H61").Value = 2189.3462
$ws.Range("I61").Value = 2076.92
$ws.Range("K61").Value = 2076.92
$ws.Range("M61").Value = -1874.92
$ws.Range("H93").Value = 13334473
$ws.Range("J93").Value = 33335092
$ws.Range("L93").Value = 33335092
$ws.Range("N93").Value = -33337588
$ws.Range("H113").Value = 2189.3462
$ws.Range("I113").Value = 2076.92
$ws.Range("K113").Value = 2076.92
$ws.Range("M113").Value = 93.07999999999993
$ws.Range("H132").Value = 3212.244
$ws.Range("I132").Value = 2645.52
$ws.Range("J132").Value = 4097.75
$ws.Range("K132").Value = 7936.559999999999
$ws.Range("L132").Value = 12293.25
$ws.Range("M132").Value = -5406.559999999999
$ws.Range("N132").Value = -17353.25
$ws.Range("H136").Value = 2726.5264
$ws.Range("I136").Value = 1950.0834
$ws.Range("K136").Value = 5850.2502
$ws.Range("M136").Value = -3300.2502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 13833.333
$ws.Range("I18").Value = 11500
$ws.Range("K18").Value = 11500
$ws.Range("M18").Value = -11327
$ws.Range("H46").Value = 57347
$ws.Range("J46").Value = 57347
$ws.Range("L46").Value = 57347
$ws.Range("N46").Value = -57809
$ws.Range("H134").Value = 57347
$ws.Range("J134").Value = 57347
$ws.Range("L134").Value = 172041
$ws.Range("N134").Value = -177111
$ws.Range("H136").Value = 18519950
$ws.Range("I136").Value = 29240830
$ws.Range("J136").Value = 2066.4546
$ws.Range("K136").Value = 87722490
$ws.Range("L136").Value = 6199.3638
$ws.Range("M136").Value = -87719940
$ws.Range("N136").Value = -11299.3638
